$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.214"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06141"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.567"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.710"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.347"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8266"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.01356"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.1597"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08204"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'0.03140"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09242"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.905"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001700"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04814"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006207"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.006286"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001098"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001503"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.722"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Value = "'0.1208"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002687"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04608"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006978"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1133"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.003252"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.01090"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006177"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.7715"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.2045"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "48CryptobidCoinCBC"
